# Bug fixes reported by Santosh last week
#
# 1) The "datetimeFigureOut" date placeholder (Slide Master + every Custom
#    Layout) was showing the stale save date "4/6/2013" -> refresh it to
#    "7/18/13".
# 2) Slide 4 ("The Team") had "Mahender Singh" and "Rishi" as two separate
#    bullets; "Rishi" doesn't belong there (he's already listed elsewhere),
#    so that line is removed and "Mahender Singh" stays as a single,
#    corrected bullet.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached date-placeholder text everywhere it appears.
# ---------------------------------------------------------------------------
$oldDate = "4/6/2013"
$newDate = "7/18/13"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Fix the "Mahender Singh" / "Rishi" bullets on slide 4.
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $shp = $slide4.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    if ($tr.Text -notlike "*Mahender Singh*Rishi*") { continue }

    # Find the two paragraphs involved.
    $paraCount = $tr.Paragraphs().Count
    $mahenderIdx = -1
    $rishiIdx = -1
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $para = $tr.Paragraphs($pi, 1)
        $txt = $para.Text.TrimEnd([char]13, [char]10)
        if ($txt -eq "Mahender Singh") { $mahenderIdx = $pi }
        if ($txt -eq "Rishi") { $rishiIdx = $pi }
    }

    if ($mahenderIdx -gt 0 -and $rishiIdx -gt 0) {
        # Delete the whole "Rishi" paragraph (text + its paragraph mark).
        $rishiPara = $tr.Paragraphs($rishiIdx, 1)
        $rishiPara.Text = ""
        $rishiPara.Delete()

        # Re-split "Mahender Singh" into "Mahender" + " " + "Singh" so the
        # run structure matches a manual retype of the surname.
        $mahenderPara = $tr.Paragraphs($mahenderIdx, 1)
        $paraText = $mahenderPara.Text.TrimEnd([char]13, [char]10)
        $spaceIdx = $paraText.IndexOf(" ")
        if ($spaceIdx -gt 0) {
            $firstWord = $paraText.Substring(0, $spaceIdx)
            $secondWord = $paraText.Substring($spaceIdx + 1)
            $startPos = $mahenderPara.Start

            $firstRange = $tr.Characters($startPos, $spaceIdx)
            $firstRange.Text = $firstWord

            $spaceRange = $tr.Characters($startPos + $spaceIdx, 1)
            $spaceRange.Text = " "

            $secondRange = $tr.Characters($startPos + $spaceIdx + 1, $secondWord.Length)
            $secondRange.Text = $secondWord
        }
    }
}
